$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-410) holds the "Förändrad" date, stored as serial date 45171
# (2023-09-02). Update it to 45172 (2023-09-03) for every data row.
$ws.Range("C2:C410").Value = 45172
